$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling data / pushing all data / mean calculation
$updates = @{
    2  = -3
    4  = -4
    5  = -2
    6  = 2
    7  = -2
    8  = -1
    10 = -6
    11 = 7
    12 = 5
    14 = 5
    16 = -1
    17 = 1
    18 = -5
    19 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
